$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 663, shifting existing rows 663-701 down to 664-702
$ws.Rows.Item(663).Insert()

# Populate the newly inserted row 663 with the new record
$ws.Cells.Item(663, 1).Value = 6
$ws.Cells.Item(663, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(663, 3).Value = "Metropolitana"
$ws.Cells.Item(663, 4).Value = 44585
$ws.Cells.Item(663, 5).Value = 13
$ws.Cells.Item(663, 6).Value = 100112031
$ws.Cells.Item(663, 7).Value = "Poroto verde"
$ws.Cells.Item(663, 8).Value = "Magnum"
$ws.Cells.Item(663, 9).Value = "Primera"
$ws.Cells.Item(663, 10).Value = 410
$ws.Cells.Item(663, 11).Value = 25000
$ws.Cells.Item(663, 12).Value = 27000
$ws.Cells.Item(663, 13).Value = 25878
$ws.Cells.Item(663, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(663, 15).Value = "Región Metropolitana"
$ws.Cells.Item(663, 16).Value = 1035
$ws.Cells.Item(663, 17).Value = 25
$ws.Cells.Item(663, 18).Value = "Hortaliza"
